$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing single-row layout so shared strings get rebuilt cleanly.
$ws.Cells.Clear()

# Block 1: parameters (A1:A3)
$ws.Range("A1").Value = "parameters"
$ws.Range("A2").Value = "model"
$ws.Range("A3").Value = "elab_method"

# Block 2: scores (E1:E4)
$ws.Range("E1").Value = "scores"
$ws.Range("E2").Value = "{'response': 'build a wall', 'prompt': 'brick', 'originality': 1.0}"
$ws.Range("E3").Value = "{'response': 'paper weight', 'prompt': 'brick', 'originality': 1.2}"
$ws.Range("E4").Value = "{'response': 'weapon', 'prompt': 'brick', 'originality': 1.4}"

# Block 3: summative (I1:I3)
$ws.Range("I1").Value = "summative"
$ws.Range("I2").Value = "n_examples"
$ws.Range("I3").Value = "originality"

# Block 4: version (M1:M2)
$ws.Range("M1").Value = "version"
$ws.Range("M2").Value = "1.3"

# Block 5: cite (Q1:Q3)
$ws.Range("Q1").Value = "cite"
$ws.Range("Q2").Value = "Organisciak, P., & Dumas, D. (2020). Open Creativity Scoring. University of Denver. https://openscoring.du.edu"
$ws.Range("Q3").Value = "Organisciak, P., Acar, S., Dumas, D., & Berthiaume, K. (2023). Beyond semantic distance: Automated scoring of divergent thinking greatly improves with large language models. Thinking Skills and Creativity, 49, 101356. https://doi.org/10.1016/j.tsc.2023.101356"
